# Atualização de bases das ligas, do dia: 26-02-2024 às 22:04
#
# The source data for "India Super League" gained a new match result and a
# brand-new fixture was inserted between the existing rows 96 and 97
# (worksheet rows), pushing the two following rows down by one. We therefore
# recompute the final state of worksheet rows 96-99 directly (rather than
# using Rows.Insert, which would fabricate an unused style entry) and write
# every cell explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 96 (id=94, match 7751752): result is now known (1-0, Home win) and
# the closing odds / P&L columns were updated; two new trailing columns
# (AB, AC) appear.
# ---------------------------------------------------------------------
$ws.Range("H96").Value = 1
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = "H"

$ws.Range("N96").Value = 2.15
$ws.Range("O96").Value = 3.25
$ws.Range("P96").Value = 3.1
$ws.Range("Q96").Value = -0.25
$ws.Range("R96").Value = 1.875
$ws.Range("S96").Value = 1.925
$ws.Range("U96").Value = 1.825
$ws.Range("V96").Value = 1.975
$ws.Range("W96").Value = 1.15
$ws.Range("X96").Value = -1
$ws.Range("Y96").Value = -1
$ws.Range("Z96").Value = 0.875
$ws.Range("AA96").Value = -1
$ws.Range("AB96").Value = -1
$ws.Range("AC96").Value = 0.9750000000000001

# ---------------------------------------------------------------------
# Row 97 (id=95): a brand-new fixture (match 7751753, Hyderabad FC vs
# Punjab FC) replaces what used to be match 7749871 in this slot - that
# older match's data is pushed down to row 98 below.
# ---------------------------------------------------------------------
$ws.Range("B97").Value = 7751753
$ws.Range("E97").Value = 45349.45833333334
$ws.Range("F97").Value = "Hyderabad FC"
$ws.Range("G97").Value = "Punjab FC"

$ws.Range("K97").Value = 4.333
$ws.Range("L97").Value = 3.6
$ws.Range("M97").Value = 1.65
$ws.Range("N97").Value = 4.333
$ws.Range("O97").Value = 3.5
$ws.Range("P97").Value = 1.65
$ws.Range("Q97").Value = 0.75
$ws.Range("R97").Value = 1.875
$ws.Range("S97").Value = 1.925

# ---------------------------------------------------------------------
# Row 98 (new, id=96): the match (7749871, Mumbai City FC vs FC Goa) that
# used to sit in row 97 - carried down unchanged, plus its odds got an
# update now that the match has more recent closing lines.
# ---------------------------------------------------------------------
$ws.Range("A98").Value = 96
$ws.Range("B98").Value = 7749871
$ws.Range("C98").Value = "India Super League"
$ws.Range("D98").Value = "India Super League"
$ws.Range("E98").Value = 45350.45833333334
$ws.Range("F98").Value = "Mumbai City FC"
$ws.Range("G98").Value = "FC Goa"

$ws.Range("K98").Value = 1.909
$ws.Range("L98").Value = 3.3
$ws.Range("M98").Value = 3.6
$ws.Range("N98").Value = 2
$ws.Range("O98").Value = 3.4
$ws.Range("P98").Value = 3.3
$ws.Range("Q98").Value = -0.5
$ws.Range("R98").Value = 1.975
$ws.Range("S98").Value = 1.825
$ws.Range("T98").Value = 2.5
$ws.Range("U98").Value = 1.9
$ws.Range("V98").Value = 1.9
$ws.Range("W98").Value = 0
$ws.Range("X98").Value = 0
$ws.Range("Y98").Value = 0
$ws.Range("Z98").Value = 0
$ws.Range("AA98").Value = 0

# Column A / E carry dedicated styles (bold+border id column, date format)
# elsewhere in the sheet - copy those formats onto the freshly written row
# instead of leaving them default-styled.
$ws.Range("A97").Copy() | Out-Null
$ws.Range("A98").PasteSpecial(-4122) | Out-Null
$ws.Range("E97").Copy() | Out-Null
$ws.Range("E98").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Row 99 (new, id=97): the match (7749767, Odisha FC vs East Bengal Club)
# that used to sit in row 98 - carried down, also with refreshed closing
# odds.
# ---------------------------------------------------------------------
$ws.Range("A99").Value = 97
$ws.Range("B99").Value = 7749767
$ws.Range("C99").Value = "India Super League"
$ws.Range("D99").Value = "India Super League"
$ws.Range("E99").Value = 45351.45833333334
$ws.Range("F99").Value = "Odisha FC"
$ws.Range("G99").Value = "East Bengal Club"

$ws.Range("K99").Value = 1.6
$ws.Range("L99").Value = 3.8
$ws.Range("M99").Value = 5
$ws.Range("N99").Value = 1.615
$ws.Range("O99").Value = 3.75
$ws.Range("P99").Value = 5
$ws.Range("Q99").Value = -0.75
$ws.Range("R99").Value = 1.8
$ws.Range("S99").Value = 2
$ws.Range("T99").Value = 2.5
$ws.Range("U99").Value = 1.8
$ws.Range("V99").Value = 2
$ws.Range("W99").Value = 0
$ws.Range("X99").Value = 0
$ws.Range("Y99").Value = 0
$ws.Range("Z99").Value = 0
$ws.Range("AA99").Value = 0

$ws.Range("A97").Copy() | Out-Null
$ws.Range("A99").PasteSpecial(-4122) | Out-Null
$ws.Range("E97").Copy() | Out-Null
$ws.Range("E99").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

Write-Host "India Super League sheet updated: rows 96-99 refreshed."
